$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-13 (M1 block unchanged text; Neutro block -> M2) ---
# Row 2
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 4.289920333333334
$ws.Range("H2").Value2 = 12.869761
$ws.Range("I2").Value2 = 0.5109350873302341
$ws.Range("J2").Value2 = 0.5109350873302341
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 260.05794
$ws.Range("N2").Value2 = 780.17382
$ws.Range("O2").Value2 = 0.9183237679872462
$ws.Range("P2").Value2 = 0.9183237679872462
$ws.Range("Q2").Value2 = 1115.62784465078
$ws.Range("R2").Value2 = 10040.65060185702
$ws.Range("S2").Value2 = 0.4692038345939933
$ws.Range("T2").Value2 = 0.4692038345939933

# Row 3
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 4.289920333333334
$ws.Range("H3").Value2 = 12.869761
$ws.Range("I3").Value2 = 0.5109350873302341
$ws.Range("J3").Value2 = 0.5109350873302341
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 3.298516
$ws.Range("N3").Value2 = 9.895548
$ws.Range("O3").Value2 = 0.011647810645144
$ws.Range("P3").Value2 = 0.011647810645144
$ws.Range("Q3").Value2 = 14.15037085822533
$ws.Range("R3").Value2 = 127.353337724028
$ws.Range("S3").Value2 = 0.00595127514918268
$ws.Range("T3").Value2 = 0.005951275149182679

# Row 4
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 4.289920333333334
$ws.Range("H4").Value2 = 12.869761
$ws.Range("I4").Value2 = 0.5109350873302341
$ws.Range("J4").Value2 = 0.5109350873302341
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 5.199381666666667
$ws.Range("N4").Value2 = 15.598145
$ws.Range("O4").Value2 = 0.01836019989751954
$ws.Range("P4").Value2 = 0.01836019989751954
$ws.Range("Q4").Value2 = 22.30493313259389
$ws.Range("R4").Value2 = 200.744398193345
$ws.Range("S4").Value2 = 0.0093808703380397
$ws.Range("T4").Value2 = 0.0093808703380397

# Row 5
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 4.289920333333334
$ws.Range("H5").Value2 = 12.869761
$ws.Range("I5").Value2 = 0.5109350873302341
$ws.Range("J5").Value2 = 0.5109350873302341
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 5.482828
$ws.Range("N5").Value2 = 16.448484
$ws.Range("O5").Value2 = 0.01936111340490499
$ws.Range("P5").Value2 = 0.01936111340490499
$ws.Range("Q5").Value2 = 23.52089532136934
$ws.Range("R5").Value2 = 211.688057892324
$ws.Range("S5").Value2 = 0.009892272168345699
$ws.Range("T5").Value2 = 0.009892272168345697

# Row 6
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 4.289920333333334
$ws.Range("H6").Value2 = 12.869761
$ws.Range("I6").Value2 = 0.5109350873302341
$ws.Range("J6").Value2 = 0.5109350873302341
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 5.580850666666667
$ws.Range("N6").Value2 = 16.742552
$ws.Range("O6").Value2 = 0.01970725374809732
$ws.Range("P6").Value2 = 0.01970725374809732
$ws.Range("Q6").Value2 = 23.94140475223022
$ws.Range("R6").Value2 = 215.472642770072
$ws.Range("S6").Value2 = 0.01006912741482319
$ws.Range("T6").Value2 = 0.01006912741482319

# Row 7
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 4.289920333333334
$ws.Range("H7").Value2 = 12.869761
$ws.Range("I7").Value2 = 0.5109350873302341
$ws.Range("J7").Value2 = 0.5109350873302341
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 3.568123
$ws.Range("N7").Value2 = 10.704369
$ws.Range("O7").Value2 = 0.01259985431708779
$ws.Range("P7").Value2 = 0.01259985431708779
$ws.Range("Q7").Value2 = 15.30696340953433
$ws.Range("R7").Value2 = 137.762670685809
$ws.Range("S7").Value2 = 0.006437707665849477
$ws.Range("T7").Value2 = 0.006437707665849476

# Row 8
$ws.Range("A8").Value2 = "M2"
$ws.Range("E8").Value2 = 2
$ws.Range("F8").Value2 = 0.6666666666666666
$ws.Range("G8").Value2 = 1.565458
$ws.Range("H8").Value2 = 4.696374
$ws.Range("I8").Value2 = 0.1864480824333443
$ws.Range("J8").Value2 = 0.1864480824333443
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 260.05794
$ws.Range("N8").Value2 = 780.17382
$ws.Range("O8").Value2 = 0.9183237679872462
$ws.Range("P8").Value2 = 0.9183237679872462
$ws.Range("Q8").Value2 = 407.1097826365199
$ws.Range("R8").Value2 = 3663.98804372868
$ws.Range("S8").Value2 = 0.1712197055941855
$ws.Range("T8").Value2 = 0.1712197055941855

# Row 9
$ws.Range("A9").Value2 = "M2"
$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = 0.6666666666666666
$ws.Range("G9").Value2 = 1.565458
$ws.Range("H9").Value2 = 4.696374
$ws.Range("I9").Value2 = 0.1864480824333443
$ws.Range("J9").Value2 = 0.1864480824333443
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 3.298516
$ws.Range("N9").Value2 = 9.895548
$ws.Range("O9").Value2 = 0.011647810645144
$ws.Range("P9").Value2 = 0.011647810645144
$ws.Range("Q9").Value2 = 5.163688260327999
$ws.Range("R9").Value2 = 46.47319434295199
$ws.Range("S9").Value2 = 0.002171711959333794
$ws.Range("T9").Value2 = 0.002171711959333794

# Row 10
$ws.Range("A10").Value2 = "M2"
$ws.Range("E10").Value2 = 2
$ws.Range("F10").Value2 = 0.6666666666666666
$ws.Range("G10").Value2 = 1.565458
$ws.Range("H10").Value2 = 4.696374
$ws.Range("I10").Value2 = 0.1864480824333443
$ws.Range("J10").Value2 = 0.1864480824333443
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 5.199381666666667
$ws.Range("N10").Value2 = 15.598145
$ws.Range("O10").Value2 = 0.01836019989751954
$ws.Range("P10").Value2 = 0.01836019989751954
$ws.Range("Q10").Value2 = 8.139413625136665
$ws.Range("R10").Value2 = 73.25472262622999
$ws.Range("S10").Value2 = 0.003423224063985403
$ws.Range("T10").Value2 = 0.003423224063985403

# Row 11
$ws.Range("A11").Value2 = "M2"
$ws.Range("E11").Value2 = 2
$ws.Range("F11").Value2 = 0.6666666666666666
$ws.Range("G11").Value2 = 1.565458
$ws.Range("H11").Value2 = 4.696374
$ws.Range("I11").Value2 = 0.1864480824333443
$ws.Range("J11").Value2 = 0.1864480824333443
$ws.Range("K11").Value2 = 3
$ws.Range("M11").Value2 = 5.482828
$ws.Range("N11").Value2 = 16.448484
$ws.Range("O11").Value2 = 0.01936111340490499
$ws.Range("P11").Value2 = 0.01936111340490499
$ws.Range("Q11").Value2 = 8.583136955223999
$ws.Range("R11").Value2 = 77.248232597016
$ws.Range("S11").Value2 = 0.003609842468119054
$ws.Range("T11").Value2 = 0.003609842468119054

# Row 12
$ws.Range("A12").Value2 = "M2"
$ws.Range("E12").Value2 = 2
$ws.Range("F12").Value2 = 0.6666666666666666
$ws.Range("G12").Value2 = 1.565458
$ws.Range("H12").Value2 = 4.696374
$ws.Range("I12").Value2 = 0.1864480824333443
$ws.Range("J12").Value2 = 0.1864480824333443
$ws.Range("K12").Value2 = 3
$ws.Range("M12").Value2 = 5.580850666666667
$ws.Range("N12").Value2 = 16.742552
$ws.Range("O12").Value2 = 0.01970725374809732
$ws.Range("P12").Value2 = 0.01970725374809732
$ws.Range("Q12").Value2 = 8.736587322938666
$ws.Range("R12").Value2 = 78.62928590644799
$ws.Range("S12").Value2 = 0.003674379671360084
$ws.Range("T12").Value2 = 0.003674379671360084

# Row 13
$ws.Range("A13").Value2 = "M2"
$ws.Range("E13").Value2 = 2
$ws.Range("F13").Value2 = 0.6666666666666666
$ws.Range("G13").Value2 = 1.565458
$ws.Range("H13").Value2 = 4.696374
$ws.Range("I13").Value2 = 0.1864480824333443
$ws.Range("J13").Value2 = 0.1864480824333443
$ws.Range("K13").Value2 = 3
$ws.Range("M13").Value2 = 3.568123
$ws.Range("N13").Value2 = 10.704369
$ws.Range("O13").Value2 = 0.01259985431708779
$ws.Range("P13").Value2 = 0.01259985431708779
$ws.Range("Q13").Value2 = 5.585746695334
$ws.Range("R13").Value2 = 50.27172025800599
$ws.Range("S13").Value2 = 0.002349218676360514
$ws.Range("T13").Value2 = 0.002349218676360514

# --- Append new rows 14-19 (Neutro sending-cluster block) ---
# Row 14
$ws.Range("A14").Value2 = "Neutro"
$ws.Range("B14").Value2 = "Cd177"
$ws.Range("C14").Value2 = "Pecam1"
$ws.Range("D14").Value2 = "ECs"
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 2.540835666666667
$ws.Range("H14").Value2 = 7.622507
$ws.Range("I14").Value2 = 0.3026168302364216
$ws.Range("J14").Value2 = 0.3026168302364217
$ws.Range("K14").Value2 = 3
$ws.Range("M14").Value2 = 260.05794
$ws.Range("N14").Value2 = 780.17382
$ws.Range("O14").Value2 = 0.9183237679872462
$ws.Range("P14").Value2 = 0.9183237679872462
$ws.Range("Q14").Value2 = 660.76448935186
$ws.Range("R14").Value2 = 5946.88040416674
$ws.Range("S14").Value2 = 0.2779002277990675
$ws.Range("T14").Value2 = 0.2779002277990676
$ws.Range("L14").Value2 = 1

# Row 15
$ws.Range("A15").Value2 = "Neutro"
$ws.Range("B15").Value2 = "Cd177"
$ws.Range("C15").Value2 = "Pecam1"
$ws.Range("D15").Value2 = "FAPs"
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 2.540835666666667
$ws.Range("H15").Value2 = 7.622507
$ws.Range("I15").Value2 = 0.3026168302364216
$ws.Range("J15").Value2 = 0.3026168302364217
$ws.Range("K15").Value2 = 3
$ws.Range("M15").Value2 = 3.298516
$ws.Range("N15").Value2 = 9.895548
$ws.Range("O15").Value2 = 0.011647810645144
$ws.Range("P15").Value2 = 0.011647810645144
$ws.Range("Q15").Value2 = 8.380987099870666
$ws.Range("R15").Value2 = 75.42888389883599
$ws.Range("S15").Value2 = 0.003524823536627527
$ws.Range("T15").Value2 = 0.003524823536627527
$ws.Range("L15").Value2 = 1

# Row 16
$ws.Range("A16").Value2 = "Neutro"
$ws.Range("B16").Value2 = "Cd177"
$ws.Range("C16").Value2 = "Pecam1"
$ws.Range("D16").Value2 = "M1"
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 2.540835666666667
$ws.Range("H16").Value2 = 7.622507
$ws.Range("I16").Value2 = 0.3026168302364216
$ws.Range("J16").Value2 = 0.3026168302364217
$ws.Range("K16").Value2 = 3
$ws.Range("M16").Value2 = 5.199381666666667
$ws.Range("N16").Value2 = 15.598145
$ws.Range("O16").Value2 = 0.01836019989751954
$ws.Range("P16").Value2 = 0.01836019989751954
$ws.Range("Q16").Value2 = 13.21077438327944
$ws.Range("R16").Value2 = 118.896969449515
$ws.Range("S16").Value2 = 0.005556105495494436
$ws.Range("T16").Value2 = 0.005556105495494437
$ws.Range("L16").Value2 = 1

# Row 17
$ws.Range("A17").Value2 = "Neutro"
$ws.Range("B17").Value2 = "Cd177"
$ws.Range("C17").Value2 = "Pecam1"
$ws.Range("D17").Value2 = "M2"
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 2.540835666666667
$ws.Range("H17").Value2 = 7.622507
$ws.Range("I17").Value2 = 0.3026168302364216
$ws.Range("J17").Value2 = 0.3026168302364217
$ws.Range("K17").Value2 = 3
$ws.Range("M17").Value2 = 5.482828
$ws.Range("N17").Value2 = 16.448484
$ws.Range("O17").Value2 = 0.01936111340490499
$ws.Range("P17").Value2 = 0.01936111340490499
$ws.Range("Q17").Value2 = 13.93096493659867
$ws.Range("R17").Value2 = 125.378684429388
$ws.Range("S17").Value2 = 0.005858998768440241
$ws.Range("T17").Value2 = 0.005858998768440242
$ws.Range("L17").Value2 = 1

# Row 18
$ws.Range("A18").Value2 = "Neutro"
$ws.Range("B18").Value2 = "Cd177"
$ws.Range("C18").Value2 = "Pecam1"
$ws.Range("D18").Value2 = "Neutro"
$ws.Range("E18").Value2 = 3
$ws.Range("F18").Value2 = 1
$ws.Range("G18").Value2 = 2.540835666666667
$ws.Range("H18").Value2 = 7.622507
$ws.Range("I18").Value2 = 0.3026168302364216
$ws.Range("J18").Value2 = 0.3026168302364217
$ws.Range("K18").Value2 = 3
$ws.Range("M18").Value2 = 5.580850666666667
$ws.Range("N18").Value2 = 16.742552
$ws.Range("O18").Value2 = 0.01970725374809732
$ws.Range("P18").Value2 = 0.01970725374809732
$ws.Range("Q18").Value2 = 14.18002442420711
$ws.Range("R18").Value2 = 127.620219817864
$ws.Range("S18").Value2 = 0.005963746661914051
$ws.Range("T18").Value2 = 0.005963746661914051
$ws.Range("L18").Value2 = 1

# Row 19
$ws.Range("A19").Value2 = "Neutro"
$ws.Range("B19").Value2 = "Cd177"
$ws.Range("C19").Value2 = "Pecam1"
$ws.Range("D19").Value2 = "sCs"
$ws.Range("E19").Value2 = 3
$ws.Range("F19").Value2 = 1
$ws.Range("G19").Value2 = 2.540835666666667
$ws.Range("H19").Value2 = 7.622507
$ws.Range("I19").Value2 = 0.3026168302364216
$ws.Range("J19").Value2 = 0.3026168302364217
$ws.Range("K19").Value2 = 3
$ws.Range("M19").Value2 = 3.568123
$ws.Range("N19").Value2 = 10.704369
$ws.Range("O19").Value2 = 0.01259985431708779
$ws.Range("P19").Value2 = 0.01259985431708779
$ws.Range("Q19").Value2 = 9.066014181453667
$ws.Range("R19").Value2 = 81.59412763308299
$ws.Range("S19").Value2 = 0.0038129279748778
$ws.Range("T19").Value2 = 0.0038129279748778
$ws.Range("L19").Value2 = 1

Write-Output "Edit complete"